$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 130
$ws.Cells.Item(130, 8).Value = 43277.145  # H130: was 38436
$ws.Cells.Item(130, 10).Value = 43277.145  # J130: was 38436
$ws.Cells.Item(130, 12).Value = 43277.145  # L130: was 38436
$ws.Cells.Item(130, 14).Value = -53317.145  # N130: was -48476

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Cells.Item(2, 8).Value = 2896.158  # H2: was 1692
$ws.Cells.Item(2, 9).Value = 2025  # I2: was 1011.129
$ws.Cells.Item(2, 10).Value = 3529.7273  # J2: was 3802.7
$ws.Cells.Item(2, 11).Value = 2025  # K2: was 1011.129
$ws.Cells.Item(2, 12).Value = 3529.7273  # L2: was 3802.7
$ws.Cells.Item(2, 13).Value = -1912  # M2: was -898.129
$ws.Cells.Item(2, 14).Value = -3755.7273  # N2: was -4028.7
# Row 32
$ws.Cells.Item(32, 8).Value = 22304.871  # H32: was 3905.81
$ws.Cells.Item(32, 9).Value = 6509.315  # I32: was 3694.704
$ws.Cells.Item(32, 10).Value = 118394.5  # J32: was 14250
$ws.Cells.Item(32, 11).Value = 6509.315  # K32: was 3694.704
$ws.Cells.Item(32, 12).Value = 118394.5  # L32: was 14250
$ws.Cells.Item(32, 13).Value = -6222.315  # M32: was -3407.704
$ws.Cells.Item(32, 14).Value = -118968.5  # N32: was -14824
# Row 80
$ws.Cells.Item(80, 8).Value = 37055  # H80: was 37105
$ws.Cells.Item(80, 10).Value = 37055  # J80: was 37105
$ws.Cells.Item(80, 12).Value = 37055  # L80: was 37105
$ws.Cells.Item(80, 14).Value = -39051  # N80: was -39101
# Row 83
$ws.Cells.Item(83, 8).Value = 37055  # H83: was 37105
$ws.Cells.Item(83, 10).Value = 37055  # J83: was 37105
$ws.Cells.Item(83, 12).Value = 111165  # L83: was 111315
$ws.Cells.Item(83, 14).Value = -121149  # N83: was -121299
# Row 116
$ws.Cells.Item(116, 8).Value = 2896.158  # H116: was 1692
$ws.Cells.Item(116, 9).Value = 2025  # I116: was 1011.129
$ws.Cells.Item(116, 10).Value = 3529.7273  # J116: was 3802.7
$ws.Cells.Item(116, 11).Value = 2025  # K116: was 1011.129
$ws.Cells.Item(116, 12).Value = 3529.7273  # L116: was 3802.7
$ws.Cells.Item(116, 13).Value = 269  # M116: was 1282.871
$ws.Cells.Item(116, 14).Value = -8117.7273  # N116: was -8390.700000000001
# Row 125
$ws.Cells.Item(125, 8).Value = 180037980  # H125: was 180032000
$ws.Cells.Item(125, 10).Value = 180037980  # J125: was 180032000
$ws.Cells.Item(125, 12).Value = 180037980  # L125: was 180032000
$ws.Cells.Item(125, 14).Value = -180047820  # N125: was -180041840
# Row 128
$ws.Cells.Item(128, 8).Value = 54980  # H128: was 0
$ws.Cells.Item(128, 10).Value = 54980  # J128: was 0
$ws.Cells.Item(128, 12).Value = 54980  # L128: was 0
$ws.Cells.Item(128, 14).Value = -64940  # N128: was ADDED
# Row 134
$ws.Cells.Item(134, 8).Value = 38645.715  # H134: was 39350
$ws.Cells.Item(134, 10).Value = 38645.715  # J134: was 39350
$ws.Cells.Item(134, 12).Value = 38645.715  # L134: was 39350
$ws.Cells.Item(134, 14).Value = -48785.715  # N134: was -49490
# Row 135
$ws.Cells.Item(135, 8).Value = 32592.834  # H135: was 34799.5
$ws.Cells.Item(135, 10).Value = 32592.834  # J135: was 34799.5
$ws.Cells.Item(135, 12).Value = 32592.834  # L135: was 34799.5
$ws.Cells.Item(135, 14).Value = -42732.834  # N135: was -44939.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Cells.Item(3, 8).Value = 2896.158  # H3: was 1692
$ws.Cells.Item(3, 9).Value = 2025  # I3: was 1011.129
$ws.Cells.Item(3, 10).Value = 3529.7273  # J3: was 3802.7
$ws.Cells.Item(3, 11).Value = 2025  # K3: was 1011.129
$ws.Cells.Item(3, 12).Value = 3529.7273  # L3: was 3802.7
$ws.Cells.Item(3, 13).Value = -1911  # M3: was -897.129
$ws.Cells.Item(3, 14).Value = -3757.7273  # N3: was -4030.7
# Row 82
$ws.Cells.Item(82, 8).Value = 22409.625  # H82: was 22415.875
$ws.Cells.Item(82, 10).Value = 29670.3  # J82: was 29680.3
$ws.Cells.Item(82, 12).Value = 29670.3  # L82: was 29680.3
$ws.Cells.Item(82, 14).Value = -30436.3  # N82: was -30446.3
# Row 85
$ws.Cells.Item(85, 8).Value = 22409.625  # H85: was 22415.875
$ws.Cells.Item(85, 10).Value = 29670.3  # J85: was 29680.3
$ws.Cells.Item(85, 12).Value = 29670.3  # L85: was 29680.3
$ws.Cells.Item(85, 14).Value = -32322.3  # N85: was -32332.3
# Row 113
$ws.Cells.Item(113, 8).Value = 4470  # H113: was 4940
$ws.Cells.Item(113, 9).Value = 4470  # I113: was 4940
$ws.Cells.Item(113, 11).Value = 4470  # K113: was 4940
$ws.Cells.Item(113, 13).Value = -2300  # M113: was -2770
# Row 122
$ws.Cells.Item(122, 8).Value = 47240  # H122: was 45331.25
$ws.Cells.Item(122, 10).Value = 47240  # J122: was 45331.25
$ws.Cells.Item(122, 12).Value = 47240  # L122: was 45331.25
$ws.Cells.Item(122, 14).Value = -57040  # N122: was -55131.25
# Row 124
$ws.Cells.Item(124, 8).Value = 42397.5  # H124: was 42450
$ws.Cells.Item(124, 10).Value = 42397.5  # J124: was 42450
$ws.Cells.Item(124, 12).Value = 42397.5  # L124: was 42450
$ws.Cells.Item(124, 14).Value = -52217.5  # N124: was -52270

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 60
$ws.Cells.Item(60, 8).Value = 23219.615  # H60: was 23218.846
$ws.Cells.Item(60, 10).Value = 25146.834  # J60: was 25146
$ws.Cells.Item(60, 12).Value = 25146.834  # L60: was 25146
$ws.Cells.Item(60, 14).Value = -26168.834  # N60: was -26168
# Row 107
$ws.Cells.Item(107, 8).Value = 916.2273  # H107: was 954.0526
$ws.Cells.Item(107, 9).Value = 950.82355  # I107: was 1058.8182
$ws.Cells.Item(107, 10).Value = 798.6  # J107: was 810
$ws.Cells.Item(107, 11).Value = 950.82355  # K107: was 1058.8182
$ws.Cells.Item(107, 12).Value = 798.6  # L107: was 810
$ws.Cells.Item(107, 13).Value = 969.17645  # M107: was 861.1818000000001
$ws.Cells.Item(107, 14).Value = -4638.6  # N107: was -4650
# Row 122
$ws.Cells.Item(122, 8).Value = 2942.2222  # H122: was 2394.7368
$ws.Cells.Item(122, 9).Value = 1750  # I122: was 1563.6364
$ws.Cells.Item(122, 10).Value = 3896  # J122: was 3537.5
$ws.Cells.Item(122, 11).Value = 5250  # K122: was 4690.9092
$ws.Cells.Item(122, 12).Value = 11688  # L122: was 10612.5
$ws.Cells.Item(122, 13).Value = -2800  # M122: was -2240.9092
$ws.Cells.Item(122, 14).Value = -16588  # N122: was -15512.5
# Row 127
$ws.Cells.Item(127, 8).Value = 56000  # H127: was 56500
$ws.Cells.Item(127, 10).Value = 56000  # J127: was 56500
$ws.Cells.Item(127, 12).Value = 56000  # L127: was 56500
$ws.Cells.Item(127, 14).Value = -65920  # N127: was -66420
# Row 130
$ws.Cells.Item(130, 8).Value = 54180  # H130: was 54125
$ws.Cells.Item(130, 10).Value = 54180  # J130: was 54125
$ws.Cells.Item(130, 12).Value = 54180  # L130: was 54125
$ws.Cells.Item(130, 14).Value = -64220  # N130: was -64165

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Cells.Item(12, 8).Value = 394.64517  # H12: was 407.46667
$ws.Cells.Item(12, 10).Value = 550.619  # J12: was 577.65
$ws.Cells.Item(12, 12).Value = 1651.857  # L12: was 1732.95
$ws.Cells.Item(12, 14).Value = -1997.857  # N12: was -2078.95
# Row 92
$ws.Cells.Item(92, 8).Value = 978.0909  # H92: was 1060
$ws.Cells.Item(92, 9).Value = 757.375  # I92: was 780
$ws.Cells.Item(92, 10).Value = 1566.6666  # J92: was 1550
$ws.Cells.Item(92, 11).Value = 2272.125  # K92: was 2340
$ws.Cells.Item(92, 12).Value = 4699.9998  # L92: was 4650
$ws.Cells.Item(92, 13).Value = -1024.125  # M92: was -1092
$ws.Cells.Item(92, 14).Value = -7195.9998  # N92: was -7146
# Row 101
$ws.Cells.Item(101, 8).Value = 12333.333  # H101: was 5914.2856
$ws.Cells.Item(101, 10).Value = 12333.333  # J101: was 5914.2856
$ws.Cells.Item(101, 12).Value = 36999.999  # L101: was 17742.8568
$ws.Cells.Item(101, 14).Value = -41867.999  # N101: was -22610.8568
# Row 131
$ws.Cells.Item(131, 8).Value = 921.0700000000001  # H131: was 917.33685
$ws.Cells.Item(131, 10).Value = 961.6882000000001  # J131: was 959.9659
$ws.Cells.Item(131, 12).Value = 2885.0646  # L131: was 2879.8977
$ws.Cells.Item(131, 14).Value = -12965.0646  # N131: was -12959.8977

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 43
$ws.Cells.Item(43, 8).Value = 3812.8462  # H43: was 3427.3845
$ws.Cells.Item(43, 10).Value = 7211.1665  # J43: was 6376
$ws.Cells.Item(43, 12).Value = 7211.1665  # L43: was 6376
$ws.Cells.Item(43, 14).Value = -7513.1665  # N43: was -6678
# Row 93
$ws.Cells.Item(93, 8).Value = 9026.1  # H93: was 12197.5
$ws.Cells.Item(93, 10).Value = 9026.1  # J93: was 12197.5
$ws.Cells.Item(93, 12).Value = 9026.1  # L93: was 12197.5
$ws.Cells.Item(93, 14).Value = -12770.1  # N93: was -15941.5
# Row 113
$ws.Cells.Item(113, 8).Value = 1683.1923  # H113: was 1807.0869
$ws.Cells.Item(113, 9).Value = 1561.579  # I113: was 1657.0588
$ws.Cells.Item(113, 10).Value = 2013.2858  # J113: was 2232.1667
$ws.Cells.Item(113, 11).Value = 1561.579  # K113: was 1657.0588
$ws.Cells.Item(113, 12).Value = 2013.2858  # L113: was 2232.1667
$ws.Cells.Item(113, 13).Value = 608.421  # M113: was 512.9412
$ws.Cells.Item(113, 14).Value = -6353.2858  # N113: was -6572.1667
# Row 122
$ws.Cells.Item(122, 8).Value = 1792.6666  # H122: was 1718.5555
$ws.Cells.Item(122, 9).Value = 1704.4348  # I122: was 1648.04
$ws.Cells.Item(122, 10).Value = 2300  # J122: was 2600
$ws.Cells.Item(122, 11).Value = 5113.3044  # K122: was 4944.12
$ws.Cells.Item(122, 12).Value = 6900  # L122: was 7800
$ws.Cells.Item(122, 13).Value = -2663.3044  # M122: was -2494.12
$ws.Cells.Item(122, 14).Value = -11800  # N122: was -12700
# Row 126
$ws.Cells.Item(126, 8).Value = 12936  # H126: was 13799.782
$ws.Cells.Item(126, 9).Value = 2907.7778  # I126: was 3084.889
$ws.Cells.Item(126, 10).Value = 18576.875  # J126: was 20687.928
$ws.Cells.Item(126, 11).Value = 8723.3334  # K126: was 9254.667000000001
$ws.Cells.Item(126, 12).Value = 55730.625  # L126: was 62063.784
$ws.Cells.Item(126, 13).Value = -6253.3334  # M126: was -6784.667000000001
$ws.Cells.Item(126, 14).Value = -60670.625  # N126: was -67003.784
# Row 128
$ws.Cells.Item(128, 8).Value = 54860  # H128: was 45614.285
$ws.Cells.Item(128, 10).Value = 54860  # J128: was 45614.285
$ws.Cells.Item(128, 12).Value = 54860  # L128: was 45614.285
$ws.Cells.Item(128, 14).Value = -64820  # N128: was -55574.285
# Row 133
$ws.Cells.Item(133, 8).Value = 38750  # H133: was 26965.857
$ws.Cells.Item(133, 10).Value = 38750  # J133: was 26965.857
$ws.Cells.Item(133, 12).Value = 38750  # L133: was 26965.857
$ws.Cells.Item(133, 14).Value = -48870  # N133: was -37085.857
# Row 135
$ws.Cells.Item(135, 8).Value = 61492.31  # H135: was 58286.47
$ws.Cells.Item(135, 10).Value = 61492.31  # J135: was 58286.47
$ws.Cells.Item(135, 12).Value = 61492.31  # L135: was 58286.47
$ws.Cells.Item(135, 14).Value = -71632.31  # N135: was -68426.47
# Row 136
$ws.Cells.Item(136, 8).Value = 17514  # H136: was 18016.838
$ws.Cells.Item(136, 10).Value = 17514  # J136: was 18016.838
$ws.Cells.Item(136, 12).Value = 52542  # L136: was 54050.514
$ws.Cells.Item(136, 14).Value = -57642  # N136: was -59150.514

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 122
$ws.Cells.Item(122, 8).Value = 4340.1113  # H122: was 4029
$ws.Cells.Item(122, 9).Value = 3792.9092  # I122: was 4267.1113
$ws.Cells.Item(122, 10).Value = 5200  # J122: was 3850.4167
$ws.Cells.Item(122, 11).Value = 11378.7276  # K122: was 12801.3339
$ws.Cells.Item(122, 12).Value = 15600  # L122: was 11551.2501
$ws.Cells.Item(122, 13).Value = -8928.7276  # M122: was -10351.3339
$ws.Cells.Item(122, 14).Value = -20500  # N122: was -16451.2501
# Row 123
$ws.Cells.Item(123, 8).Value = 47492  # H123: was 42500
$ws.Cells.Item(123, 10).Value = 47492  # J123: was 42500
$ws.Cells.Item(123, 12).Value = 47492  # L123: was 42500
$ws.Cells.Item(123, 14).Value = -57292  # N123: was -52300
# Row 125
$ws.Cells.Item(125, 8).Value = 49891.668  # H125: was 40357.5
$ws.Cells.Item(125, 10).Value = 49891.668  # J125: was 40357.5
$ws.Cells.Item(125, 12).Value = 49891.668  # L125: was 40357.5
$ws.Cells.Item(125, 14).Value = -59731.668  # N125: was -50197.5
# Row 134
$ws.Cells.Item(134, 8).Value = 70214.5  # H134: was 51606.332
$ws.Cells.Item(134, 9).Value = 0  # I134: was 14390
$ws.Cells.Item(134, 11).Value = 0  # K134: was 14390
$ws.Cells.Item(134, 13).ClearContents()  # M134: was -9320

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 104
$ws.Cells.Item(104, 8).Value = 27833.334  # H104: was 30725
$ws.Cells.Item(104, 10).Value = 27833.334  # J104: was 30725
$ws.Cells.Item(104, 12).Value = 27833.334  # L104: was 30725
$ws.Cells.Item(104, 14).Value = -34821.334  # N104: was -37713
# Row 122
$ws.Cells.Item(122, 8).Value = 13466.723  # H122: was 10231.5
$ws.Cells.Item(122, 9).Value = 16915.54  # I122: was 12358.444
$ws.Cells.Item(122, 10).Value = 4499.8  # J122: was 3850.6667
$ws.Cells.Item(122, 11).Value = 50746.62  # K122: was 37075.33199999999
$ws.Cells.Item(122, 12).Value = 13499.4  # L122: was 11552.0001
$ws.Cells.Item(122, 13).Value = -48296.62  # M122: was -34625.33199999999
$ws.Cells.Item(122, 14).Value = -18399.4  # N122: was -16452.0001
# Row 129
$ws.Cells.Item(129, 8).Value = 39322.25  # H129: was 39329
$ws.Cells.Item(129, 10).Value = 39322.25  # J129: was 39329
$ws.Cells.Item(129, 12).Value = 39322.25  # L129: was 39329
$ws.Cells.Item(129, 14).Value = -49322.25  # N129: was -49329
